$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) -- worksheet 1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 2007
$ws1.Range("F5").Value  = 269
$ws1.Range("F7").Value  = 3090
$ws1.Range("F13").Value = 155
$ws1.Range("F14").Value = 155
$ws1.Range("F15").Value = 10188
$ws1.Range("F18").Value = 12
$ws1.Range("F20").Value = 8086
$ws1.Range("F21").Value = 12692
$ws1.Range("F24").Value = 25
$ws1.Range("F27").Value = 600
$ws1.Range("F29").Value = 420
$ws1.Range("F30").Value = 2830
$ws1.Range("F32").Value = 239
$ws1.Range("F33").Value = 7986
$ws1.Range("F34").Value = 1571
$ws1.Range("F36").Value = 71
$ws1.Range("F39").Value = 1447
$ws1.Range("F43").Value = 641

# Sheet "演出" (Performances) -- worksheet 2
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F13").Value = 72
$ws2.Range("F16").Value = 102
$ws2.Range("F17").Value = 18

# Sheet "全部类型" (All Types) -- worksheet 4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 2007
$ws4.Range("F8").Value  = 269
$ws4.Range("F10").Value = 3090
$ws4.Range("F16").Value = 155
$ws4.Range("F17").Value = 155
$ws4.Range("F18").Value = 10188
$ws4.Range("F20").Value = 12
$ws4.Range("F22").Value = 8086
$ws4.Range("F23").Value = 12692
$ws4.Range("F25").Value = 25
$ws4.Range("F28").Value = 600
$ws4.Range("F31").Value = 2830
$ws4.Range("F35").Value = 239
$ws4.Range("F36").Value = 7986
$ws4.Range("F38").Value = 71
$ws4.Range("F42").Value = 102
$ws4.Range("F43").Value = 18
$ws4.Range("F47").Value = 641
